# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the bold/centered/bordered header style from row 1 (fonts/borders/cellXfs
#    collapse back down to a single default style) and clear the "Unnamed: 0" label
#    that used to sit in A1.
$ws.Range("A1:AO1").ClearFormats()
$ws.Range("A1").ClearContents()

# 2. Corrected numeric values across the Revisit count / Fixation count / Dwell time
#    (ms) / Dwell time (%) / Fixation duration (ms) rows.
$changes = @{
    "G3" = 8;    "H3" = 19;   "J3" = 4;    "M3" = 5;    "O3" = 20;
    "P3" = 15;   "Q3" = 12;   "R3" = 9;    "S3" = 6;    "T3" = 5;
    "V3" = 0;    "AB3" = 18;  "AC3" = 5;   "AK3" = 2;   "AM3" = 2;

    "G4" = 34;   "H4" = 158;  "J4" = 6;    "M4" = 8;    "O4" = 130;
    "P4" = 88;   "Q4" = 31;   "R4" = 16;   "S4" = 15;   "T4" = 12;
    "V4" = 1;    "AB4" = 104; "AC4" = 7;   "AK4" = 5;   "AM4" = 3;

    "G5" = 9677.690000000001; "H5" = 44847.86; "J5" = 3136.66; "M5" = 3421.25;
    "O5" = 37256.45; "P5" = 27562.87; "Q5" = 11311.47; "R5" = 5723.09;
    "S5" = 6790.13;  "T5" = 3220.09;  "V5" = 467.2;    "AB5" = 44530.27;
    "AC5" = 1751.99; "AK5" = 2270.03; "AM5" = 1184.53;

    "E6" = 1.48;  "G6" = 7.22;  "H6" = 33.48; "I6" = 2.45;  "J6" = 2.34;
    "K6" = 0.85;  "M6" = 2.55;  "N6" = 2.02;  "O6" = 27.81; "P6" = 20.57;
    "Q6" = 8.44;  "R6" = 4.27;  "S6" = 5.07;  "T6" = 2.4;   "V6" = 0.35;
    "W6" = 0.54;  "X6" = 0.67;  "Y6" = 0.24;  "Z6" = 1.28;  "AA6" = 1.71;
    "AB6" = 33.24; "AC6" = 1.31; "AD6" = 0.49; "AE6" = 0.7; "AF6" = 1.91;
    "AG6" = 1.03; "AK6" = 1.69; "AL6" = 0.65; "AM6" = 0.88; "AN6" = 1.06;
    "AO6" = 0.15;

    "G7" = 284.64; "H7" = 283.85; "J7" = 522.78; "M7" = 427.66; "O7" = 286.59;
    "P7" = 313.21; "Q7" = 364.89; "R7" = 357.69; "S7" = 452.68; "T7" = 268.34;
    "V7" = 467.2;  "AB7" = 428.18; "AC7" = 250.28; "AK7" = 454.01; "AM7" = 394.84;
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

# 3. Drop the three trailing all-blank rows (10-12), shrinking the used range down
#    to A1:AO9.
$ws.Range("A10:A12").EntireRow.Delete()
